$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("func")
$ws.Activate()

# B9: new "void print_result();" entry (C++ declaration paired with A9's
# "public void print_result()"), matching the style used by its neighbours.
$ws.Range("B9").Value = "void print_result();"
$ws.Range("A9").Copy()
$ws.Range("B9").PasteSpecial(-4122)

# Recolor the "print_result" function name and the trailing "();" to match
# the other C++ declarations (e.g. B8 "void ACS_Strategy();").
$ws.Range("B9").Characters(6, 12).Font.Color = 8020480
$ws.Range("B9").Characters(18, 3).Font.Color = 526344

# B10: mirror A10 ("public void check_answer()") into B10.
$ws.Range("A10").Copy()
$ws.Range("B10").PasteSpecial(-4104)
$ws.Range("A10").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("A11").Select() | Out-Null
